$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 197.2
$ws.Range("I6").Value = 185.88889
$ws.Range("J6").Value = 299
$ws.Range("K6").Value = 557.6666700000001
$ws.Range("L6").Value = 897
$ws.Range("M6").Value = -445.6666700000001
$ws.Range("N6").Value = -1121
$ws.Range("H15").Value = 1749.9333
$ws.Range("I15").Value = 1749.9333
$ws.Range("K15").Value = 5249.7999
$ws.Range("M15").Value = -5080.7999
$ws.Range("H31").Value = 3099.8
$ws.Range("I31").Value = 3110.889
$ws.Range("K31").Value = 9332.667000000001
$ws.Range("M31").Value = -9102.667000000001
$ws.Range("H33").Value = 169.75
$ws.Range("I33").Value = 185
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = 185
$ws.Range("L33").Value = 2
$ws.Range("M33").Value = 44
$ws.Range("N33").Value = -460
$ws.Range("H40").Value = 1666.3334
$ws.Range("I40").Value = 1666.3334
$ws.Range("K40").Value = 1666.3334
$ws.Range("M40").Value = -1491.3334
$ws.Range("H64").Value = 3488.5
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 3488.5
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H135").Value = 1614.6666
$ws.Range("I135").Value = 1710.1666
$ws.Range("K135").Value = 15391.4994
$ws.Range("M135").Value = -12856.4994
$ws.Range("H138").Value = 3136.6562
$ws.Range("I138").Value = 1468.5333
$ws.Range("K138").Value = 4405.5999
$ws.Range("M138").Value = 734.4000999999998
$ws.Range("H141").Value = 2418.7693
$ws.Range("I141").Value = 2216.6667
$ws.Range("J141").Value = 4844
$ws.Range("K141").Value = 6650.000100000001
$ws.Range("L141").Value = 14532
$ws.Range("M141").Value = -1470.000100000001
$ws.Range("N141").Value = -24892

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3502999
$ws.Range("I32").Value = 3336475.2
$ws.Range("K32").Value = 3336475.2
$ws.Range("M32").Value = -3336188.2
$ws.Range("H38").Value = 93499.5
$ws.Range("I38").Value = 93499.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 93499.5
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -93032.5
$ws.Range("H74").Value = 805.4
$ws.Range("I74").Value = 779.2857
$ws.Range("K74").Value = 779.2857
$ws.Range("M74").Value = 94.71429999999998
$ws.Range("H77").Value = 805.4
$ws.Range("I77").Value = 779.2857
$ws.Range("K77").Value = 3896.4285
$ws.Range("M77").Value = 471.5715
$ws.Range("H111").Value = 60000
$ws.Range("J111").Value = 60000
$ws.Range("L111").Value = 60000
$ws.Range("N111").Value = -68180
$ws.Range("H132").Value = 1818.9
$ws.Range("I132").Value = 1774.25
$ws.Range("K132").Value = 5322.75
$ws.Range("M132").Value = -2792.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 795
$ws.Range("I99").Value = 795
$ws.Range("K99").Value = 795
$ws.Range("M99").Value = 703
$ws.Range("H107").Value = 2757.889
$ws.Range("I107").Value = 2862.2856
$ws.Range("J107").Value = 2392.5
$ws.Range("K107").Value = 2862.2856
$ws.Range("L107").Value = 2392.5
$ws.Range("M107").Value = -942.2856000000002
$ws.Range("N107").Value = -6232.5
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("N132").Value = 0
$ws.Range("H134").Value = 4710.48
$ws.Range("I134").Value = 4231.5713
$ws.Range("J134").Value = 7224.75
$ws.Range("K134").Value = 12694.7139
$ws.Range("L134").Value = 21674.25
$ws.Range("M134").Value = -10159.7139
$ws.Range("N134").Value = -26744.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 199.66667
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H42").Value = 36666.668
$ws.Range("J42").Value = 37500
$ws.Range("L42").Value = 37500
$ws.Range("N42").Value = -38686
$ws.Range("H99").Value = 1099.5
$ws.Range("I99").Value = 799
$ws.Range("K99").Value = 799
$ws.Range("M99").Value = 699
$ws.Range("H107").Value = 605
$ws.Range("I107").Value = 211
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 211
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 1709
$ws.Range("N107").Value = -4839
$ws.Range("H126").Value = 1099.5
$ws.Range("I126").Value = 799
$ws.Range("K126").Value = 2397
$ws.Range("M126").Value = 73
$ws.Range("H134").Value = 1986
$ws.Range("I134").Value = 1986
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5958
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -3423

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2044.75
$ws.Range("I46").Value = 2066.3333
$ws.Range("J46").Value = 1980
$ws.Range("K46").Value = 6198.999899999999
$ws.Range("L46").Value = 5940
$ws.Range("M46").Value = -6107.999899999999
$ws.Range("N46").Value = -6122
$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 1000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1814
$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 1000
$ws.Range("K89").Value = 9000
$ws.Range("M89").Value = -3072
$ws.Range("H107").Value = 2184.7693
$ws.Range("J107").Value = 2926
$ws.Range("L107").Value = 8778
$ws.Range("N107").Value = -12618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 32500
$ws.Range("J55").Value = 32500
$ws.Range("L55").Value = 32500
$ws.Range("N55").Value = -33154
$ws.Range("H107").Value = 752
$ws.Range("I107").Value = 514.3333
$ws.Range("J107").Value = 989.6667
$ws.Range("K107").Value = 514.3333
$ws.Range("L107").Value = 989.6667
$ws.Range("M107").Value = 1405.6667
$ws.Range("N107").Value = -4829.6667
$ws.Range("H113").Value = 5141
$ws.Range("I113").Value = 1448
$ws.Range("J113").Value = 6987.5
$ws.Range("K113").Value = 1448
$ws.Range("L113").Value = 6987.5
$ws.Range("M113").Value = 722
$ws.Range("N113").Value = -11327.5
$ws.Range("H132").Value = 3765.3635
$ws.Range("I132").Value = 3597.6
$ws.Range("K132").Value = 10792.8
$ws.Range("M132").Value = -8262.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1649.5
$ws.Range("I13").Value = 300
$ws.Range("J13").Value = 2999
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 2999
$ws.Range("M13").Value = -160
$ws.Range("N13").Value = -3279
$ws.Range("H107").Value = 1800
$ws.Range("I107").Value = 1520.4
$ws.Range("J107").Value = 2499
$ws.Range("K107").Value = 4561.200000000001
$ws.Range("L107").Value = 7497
$ws.Range("M107").Value = -2641.200000000001
$ws.Range("N107").Value = -11337
$ws.Range("H126").Value = 1617.1666
$ws.Range("I126").Value = 1272.2858
$ws.Range("K126").Value = 3816.8574
$ws.Range("M126").Value = -1346.8574
$ws.Range("H132").Value = 2023.3334
$ws.Range("I132").Value = 2023.3334
$ws.Range("K132").Value = 6070.0002
$ws.Range("M132").Value = -3540.0002
$ws.Range("H136").Value = 2288.111
$ws.Range("I136").Value = 2507.9546
$ws.Range("J136").Value = 1320.8
$ws.Range("K136").Value = 7523.8638
$ws.Range("L136").Value = 3962.4
$ws.Range("M136").Value = -4973.8638
$ws.Range("N136").Value = -9062.4
